$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (header "Förändrad") holds a date serial that was updated
# from 45181 (2023-09-12) to 45182 (2023-09-13) for every data row.
$oldValue = 45181
$newValue = 45182

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
